$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E2").Value = 102
$ws.Range("F2").Value = 74
$ws.Range("H2").Value = 80

$ws.Range("E5").Value = 142
$ws.Range("F5").Value = 98
$ws.Range("H5").Value = 109

$ws.Range("F7").Value = 22
$ws.Range("H7").Value = 27

$ws.Range("F8").Value = 4
$ws.Range("H8").Value = 7

$ws.Range("E10").Value = 608
$ws.Range("F10").Value = 311
$ws.Range("H10").Value = 407

$ws.Range("E11").Value = 393
$ws.Range("F11").Value = 220
$ws.Range("H11").Value = 284

$ws.Range("E12").Value = 604
$ws.Range("F12").Value = 340
$ws.Range("H12").Value = 426

$ws.Range("E13").Value = 144
$ws.Range("F13").Value = 79
$ws.Range("H13").Value = 113

$ws.Range("E14").Value = 132
$ws.Range("F14").Value = 74
$ws.Range("H14").Value = 108

$ws.Range("E15").Value = 182
$ws.Range("F15").Value = 78
$ws.Range("H15").Value = 128

$ws.Range("E16").Value = 212
$ws.Range("F16").Value = 111
$ws.Range("H16").Value = 159

$ws.Range("F17").Value = 60
$ws.Range("H17").Value = 84

$ws.Range("E18").Value = 54
$ws.Range("F18").Value = 27
$ws.Range("H18").Value = 44

$ws.Range("F19").Value = 7
$ws.Range("H19").Value = 10

$ws.Range("E20").Value = 91
$ws.Range("F20").Value = 33
$ws.Range("H20").Value = 70

$ws.Range("E21").Value = 143
$ws.Range("F21").Value = 82
$ws.Range("H21").Value = 113

$ws.Range("E22").Value = 178
$ws.Range("F22").Value = 98
$ws.Range("H22").Value = 140

$ws.Range("E23").Value = 209
$ws.Range("F23").Value = 105
$ws.Range("H23").Value = 156

$ws.Range("E24").Value = 230
$ws.Range("F24").Value = 127
$ws.Range("H24").Value = 157

$ws.Range("E25").Value = 293
$ws.Range("F25").Value = 155
$ws.Range("H25").Value = 215

$ws.Range("E26").Value = 165
$ws.Range("F26").Value = 102
$ws.Range("H26").Value = 127

$ws.Range("E27").Value = 348
$ws.Range("F27").Value = 184
$ws.Range("H27").Value = 265

$ws.Range("E28").Value = 209
$ws.Range("F28").Value = 95
$ws.Range("H28").Value = 147

$ws.Range("E29").Value = 175
$ws.Range("F29").Value = 105
$ws.Range("H29").Value = 146

$ws.Range("E30").Value = 226
$ws.Range("F30").Value = 135
$ws.Range("H30").Value = 187

$ws.Range("F32").Value = 118
$ws.Range("H32").Value = 156

$ws.Range("E33").Value = 309
$ws.Range("F33").Value = 164
$ws.Range("H33").Value = 253

$ws.Range("E34").Value = 230
$ws.Range("F34").Value = 156
$ws.Range("H34").Value = 195

$ws.Range("E35").Value = 161
$ws.Range("F35").Value = 103
$ws.Range("H35").Value = 130

$ws.Range("F36").Value = 49
$ws.Range("H36").Value = 59

$ws.Range("E37").Value = 173
$ws.Range("F37").Value = 90
$ws.Range("H37").Value = 127

$ws.Range("E39").Value = 187
$ws.Range("F39").Value = 96
$ws.Range("H39").Value = 147

$ws.Range("E40").Value = 278
$ws.Range("F40").Value = 131
$ws.Range("H40").Value = 211

$ws.Range("E41").Value = 410
$ws.Range("F41").Value = 203
$ws.Range("H41").Value = 295

$ws.Range("F42").Value = 228
$ws.Range("H42").Value = 289

$ws.Range("F43").Value = 68
$ws.Range("H43").Value = 95

$ws.Range("E44").Value = 325
$ws.Range("F44").Value = 169
$ws.Range("H44").Value = 237

$ws.Range("E45").Value = 161
$ws.Range("F45").Value = 84
$ws.Range("H45").Value = 123

$ws.Range("E46").Value = 344
$ws.Range("F46").Value = 192
$ws.Range("H46").Value = 255

$ws.Range("E47").Value = 484
$ws.Range("F47").Value = 258
$ws.Range("H47").Value = 350

$ws.Range("E48").Value = 232
$ws.Range("F48").Value = 101
$ws.Range("H48").Value = 145

$ws.Range("E49").Value = 301
$ws.Range("F49").Value = 139
$ws.Range("H49").Value = 226

$ws.Range("F50").Value = 131
$ws.Range("H50").Value = 202

$ws.Range("F51").Value = 116
$ws.Range("H51").Value = 190

$ws.Range("F52").Value = 13
$ws.Range("H52").Value = 21
